$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump Version 0.1.0 -> 0.2.0 ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "0.2.0"

# --- Elements sheet edits ---
$ws = $wb.Worksheets.Item("Elements")

# Row 6 (Extension.value[x]): Min 1 -> 0
$ws.Range("E6").Value = "0"

# Row 7 (Extension.value[x] / valueCodeableConcept slice): Min 1 -> 0, Must Support? -> Y
$ws.Range("E7").Value = "0"
$ws.Range("G7").Value = "Y"

# AutoFilter: Must Support? (col G, colId 6) not-equal to a single space;
# Slicing Discriminator (col AA, colId 26) blank
$space = " "
$notEqSpace = "<>" + $space
$ws.Range("A1:AJ7").AutoFilter(7, $notEqSpace)
$blanks = @("")
$ws.Range("A1:AJ7").AutoFilter(27, $blanks, 7)

# Hide detail rows 2-6 (slice summary row 7 stays visible).
# Must happen AFTER AutoFilter, since AutoFilter recomputes row visibility
# from the live filter criteria otherwise.
$ws.Range("A2:A6").EntireRow.Hidden = $true

# Persist the (normally hidden) AutoFilter defined name scoped to the Elements sheet
$fdb = $ws.Names.Add("_xlnm._FilterDatabase", "=Elements!`$A`$1:`$AJ`$7")
$fdb.Visible = $false

# Conditional formatting over the detail rows (A2:AI6)
$cfRange = $ws.Range("A2:AI6")

$cf1 = $cfRange.FormatConditions.Add(2, 0, '$G2<>"Y"')
$cf1.Interior.ColorIndex = 22

$cf2 = $cfRange.FormatConditions.Add(2, 0, '$Q2<>""')
$cf2.Font.ColorIndex = 22
$cf2.Font.Italic = $true
